# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) timestamps on the per-language
# sheets to reflect the latest handback run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 21:07:45"
$wsZhCn.Range("H2").Value = "2016-03-21 21:08:08"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 21:07:49"
$wsDeDe.Range("H2").Value = "2016-03-21 21:08:15"
